# Updated cryptos list - applying latest price/volume(1h) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.199.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.83%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.026.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.67%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.53'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.68'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0808'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.13'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.323.00'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.848'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.12'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.44'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.024.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.157.11'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0861'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.07'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.36'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.50'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.46%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.80'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.92%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0669'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +9.45%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.57'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +14.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.51'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.60'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.43%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0970'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.78%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.17%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.62'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.52'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.57'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.375.14'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.16%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.14'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +16.76%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.16'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.31%  '
